$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Scanner" to "Session"
$ws.Name = "Session"

# Remove the last logged row (row 3, A3:F3) - the used range shrinks
# from A1:F3 down to A1:F2.
$ws.Rows.Item(3).Delete()
